# ZSS-1391 test case: add a "right alignment" + "vertical center" merged-cell
# demo, next to the existing "merge + hidden" demos, on Sheet1.
#
# Layout after the edit:
#   Row 10: plain "click here" label (like A1), the visible trigger for the
#           new hidden/merged block below it.
#   Row 11 (hidden): A11:C11 merged "right alignment" label, right aligned.
#   Column H (new): width 20, hidden -> holds the "vertical center" merged
#                   label H5:H8, vertically centered.
#   Column I (new): holds a "click here" merged label I5:I8 (like G5:G8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: visible trigger label, mirrors A1 -----------------------------
$ws.Range("A10").Value = "click here"

# --- Row 11 (hidden): A11:C11 merged "right alignment" demo ---------------
$ws.Range("A11:C11").Merge()
$ws.Range("A11").Value = "right alignment"
$ws.Range("A11:C11").HorizontalAlignment = -4152   # xlRight
$ws.Rows("11:11").Hidden = $true

# --- new hidden column H, used for the "vertical center" merged label -----
$ws.Columns("H:H").ColumnWidth = 19.17
$ws.Columns("H:H").Hidden = $true

# --- H5:H8 merged cell: vertical-center alignment demo --------------------
$ws.Range("H5:H8").Merge()
$ws.Range("H5").Value = "vertical center"
$ws.Range("H5:H8").VerticalAlignment = -4108   # xlCenter

# --- I5:I8 merged cell: plain "click here" label, mirrors G5:G8 -----------
$ws.Range("I5:I8").Merge()
$ws.Range("I5").Value = "click here"

# --- refresh selection to match the new "click here" trigger column -------
$ws.Range("I5:I8").Select()
